# ---------------------------------------------------------------------------
# Applies the "Update cover and README.md" style-sheet revision:
#   * Re-center the cover-page tab stops (2520/5040 -> 2500/5000 twips).
#   * Re-point the inside-cover right tab stop (5040 -> 5000 twips).
#   * Swap the Abstract-Message East-Asian font to "楷体" (drop the
#     "_gb2312" suffix variant).
#   * Introduce a new "Abstract Underline Message" paragraph style (and
#     wire its character-style link both ways).
#   * Strip the redundant explicit justification overrides from the
#     Abstract / Keywords styles (they already inherit "both" from
#     Abstract Message).
#   * Give AbstractChar an explicit "no underline" run override, and link
#     AbstractUnderlineChar back to the new paragraph style.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$styles = $d.Styles

# --- 1. Front Cover Text: re-center tab stops -----------------------------
$frontCoverText = $styles.Item("FrontCoverText")
$frontCoverText.ParagraphFormat.TabStops.ClearAll()
$frontCoverText.ParagraphFormat.TabStops.Add(125)
$frontCoverText.ParagraphFormat.TabStops.Add(250)

# --- 2. Inside Front Cover Text: move the right tab stop ------------------
$insideFrontCoverText = $styles.Item("InsideFrontCoverText")
$insideFrontCoverText.ParagraphFormat.TabStops.ClearAll()
$insideFrontCoverText.ParagraphFormat.TabStops.Add(0.05)
$insideFrontCoverText.ParagraphFormat.TabStops.Add(250)

# --- 3. Abstract Message: normalise the East-Asian font name --------------
$abstractMessage = $styles.Item("AbstractMessage")
$abstractMessage.Font.NameFarEast = "楷体"

# --- 4. New style: Abstract Underline Message ------------------------------
$abstractUnderlineMessage = $styles.Add("AbstractUnderlineMessage", 1)
$abstractUnderlineMessage.NameLocal = "Abstract Underline Message"
$abstractUnderlineMessage.BaseStyle = "AbstractMessage"
$abstractUnderlineMessage.NextParagraphStyle = "Abstract"
$abstractUnderlineMessage.QuickStyle = $true
$abstractUnderlineMessage.Font.Underline = 1

# --- 5. Abstract: drop the redundant explicit "both" justification --------
$oldAbstract = $styles.Item("Abstract")
$oldAbstract.Delete()
$abstract = $styles.Add("Abstract", 1)
$abstract.NameLocal = "Abstract"
$abstract.BaseStyle = "AbstractMessage"
$abstract.NextParagraphStyle = "Abstract"
$abstract.QuickStyle = $true
$abstract.ParagraphFormat.CharacterUnitFirstLineIndent = 200

# --- 6. Keywords: drop the redundant explicit "both" justification --------
$oldKeywords = $styles.Item("Keywords")
$oldKeywords.Delete()
$keywords = $styles.Add("Keywords", 1)
$keywords.NameLocal = "Keywords"
$keywords.BaseStyle = "AbstractMessage"
$keywords.NextParagraphStyle = "Keywords"
$keywords.QuickStyle = $true

# --- 7. AbstractChar: explicit "no underline" run override -----------------
$abstractChar = $styles.Item("AbstractChar")
$abstractChar.Font.Underline = 0

# --- 8. AbstractUnderlineChar: link back to the new paragraph style --------
$abstractUnderlineChar = $styles.Item("AbstractUnderlineChar")
$abstractUnderlineChar.LinkStyle = "AbstractUnderlineMessage"
